$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1554434735375247
$ws.Range("C2").Value = 0.05231270169004087
$ws.Range("D2").Value = 16.98373111632243
$ws.Range("E2").Value = 0.4998867070740569
$ws.Range("G2").Value = 17.69137399862405
